# issue #5: stock data from json to db
# Adds "category", "source_file" and "index" columns to the 股票 (stock)
# sheet: a new "category" column is inserted right after
# "property_category" (pushing date/legislator_name/legislator_id one
# column to the right), and "source_file" + "index" columns are appended
# at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before the existing "date" column (column I) -
# this shifts date/legislator_name/legislator_id from I/J/K to J/K/L
# and copies the column formatting (header + data styles) along the way.
$ws.Columns("I:I").Insert()

# Append two more (empty, so far) columns at M/N for "source_file" and
# "index" - inserting (rather than just writing past the used range)
# picks up the same header/data cell styles used by the rest of the row.
$ws.Columns("M:N").Insert()

$ws.Range("I1").Value2 = "category"
$ws.Range("M1").Value2 = "source_file"
$ws.Range("N1").Value2 = "index"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $idx = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 9).Value2 = "normal"
    $ws.Cells.Item($r, 13).Value2 = "tmp22571"
    $ws.Cells.Item($r, 14).Value2 = $idx
}
